# Selumetinib SRB data - add new replicate sheet S9, fix S8's date label,
# tighten a shared-formula range on S8, and refresh sheet selections.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) S8: fix the Day-1 date label in A3 (was "8/8/19", should be "12/8/19")
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("S8")
$ws8.Range("A3").Value = "Day  1 date: 12/8/19"

# ---------------------------------------------------------------------
# 2) S8: the shared formula in C17 used to stretch over C17:K19 but only
#    row 17 actually uses it now - rewrite each cell in row 17 with its
#    own (non-shared) AVERAGE formula over the same source range.
# ---------------------------------------------------------------------
$ws8.Range("C17").Formula = "=AVERAGE(C9:C13)"
$ws8.Range("D17").Formula = "=AVERAGE(D9:D13)"
$ws8.Range("E17").Formula = "=AVERAGE(E9:E13)"
$ws8.Range("F17").Formula = "=AVERAGE(F9:F13)"
$ws8.Range("G17").Formula = "=AVERAGE(G9:G13)"
$ws8.Range("H17").Formula = "=AVERAGE(H9:H13)"
$ws8.Range("I17").Formula = "=AVERAGE(I9:I13)"
$ws8.Range("J17").Formula = "=AVERAGE(J9:J13)"
$ws8.Range("K17").Formula = "=AVERAGE(K9:K13)"

# ---------------------------------------------------------------------
# 3) S8: selection moves from E20 to G25, and it is no longer the active
#    (tab-selected) sheet once S9 is added - select its new cell now,
#    the later work on S9 will take over as the active tab.
# ---------------------------------------------------------------------
$ws8.Range("G25").Select()

# ---------------------------------------------------------------------
# 4) Add the new S9 worksheet at the end of the workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws9 = $wb.Worksheets.Add($null, $lastSheet)
$ws9.Name = "S9"

# Header block (rows 1-5) is identical in style to every other replicate
# sheet - copy the formatting straight from S8 then overwrite the text.
$ws8.Range("A1:C5").Copy()
$ws9.Range("A1").PasteSpecial(-4122)

$ws9.Range("A1").Value = "Selumetinib"
$ws9.Range("A2").Value = "Replicate # 9"
$ws9.Range("A3").Value = "Day  1 date: 22/8/19"
$ws9.Range("A4").Value = "Passage #12"

# Data block rows 7-14 and the highlighted sub-block (rows 8-13) copy the
# same conditional fills used throughout the other replicate sheets.
$ws8.Range("B9:K13").Copy()
$ws9.Range("B8").PasteSpecial(-4122)

# Row 7 - percentage confluence values (plain, no fill)
$row7 = @(0.062,0.073,0.075,0.07,0.071,0.071,0.061,0.077,0.082,0.092,0.109,0.1)
for ($i = 0; $i -lt $row7.Length; $i++) { $ws9.Cells.Item(7, $i + 1).Value = $row7[$i] }

# Row 8
$ws9.Range("A8").Value = 0.075
$row8 = @(1.359,1.2010000000000001,1.23,0.98499999999999999,0.88100000000000001,0.63100000000000001,0.64800000000000002,0.42,0.29799999999999999,0.29899999999999999)
for ($i = 0; $i -lt $row8.Length; $i++) { $ws9.Cells.Item(8, $i + 2).Value = $row8[$i] }
$ws9.Range("L8").Value = 0.113

# Row 9
$ws9.Range("A9").Value = 0.077
$row9 = @(1.3879999999999999,1.2390000000000001,1.1379999999999999,1.018,0.84899999999999998,0.63100000000000001,0.55400000000000005,0.51100000000000001,0.32400000000000001,0.33100000000000002)
for ($i = 0; $i -lt $row9.Length; $i++) { $ws9.Cells.Item(9, $i + 2).Value = $row9[$i] }
$ws9.Range("L9").Value = 0.12

# Row 10
$ws9.Range("A10").Value = 0.067
$row10 = @(1.361,1.238,1.1200000000000001,0.91900000000000004,0.86499999999999999,0.61599999999999999,0.47799999999999998,0.438,0.32,0.34799999999999998)
for ($i = 0; $i -lt $row10.Length; $i++) { $ws9.Cells.Item(10, $i + 2).Value = $row10[$i] }
$ws9.Range("L10").Value = 0.11

# Row 11 (G11 is the outlier that gets its highlight cleared below)
$ws9.Range("A11").Value = 0.086
$row11 = @(1.2150000000000001,1.1459999999999999,1.0249999999999999,0.874,0.82799999999999996,1.3160000000000001,0.52500000000000002,0.40100000000000002,0.30499999999999999,0.26900000000000002)
for ($i = 0; $i -lt $row11.Length; $i++) { $ws9.Cells.Item(11, $i + 2).Value = $row11[$i] }
$ws9.Range("L11").Value = 0.099

# Row 12
$ws9.Range("A12").Value = 0.073
$row12 = @(1.341,1.2529999999999999,0.99199999999999999,0.89100000000000001,0.85699999999999998,0.68899999999999995,0.57999999999999996,0.40500000000000003,0.30099999999999999,0.38)
for ($i = 0; $i -lt $row12.Length; $i++) { $ws9.Cells.Item(12, $i + 2).Value = $row12[$i] }
$ws9.Range("L12").Value = 0.099

# Row 13
$row13 = @(0.087,0.076,0.067,0.074,0.066,0.10100000000000001,0.161,0.123,0.094,0.083,0.08,0.086)
for ($i = 0; $i -lt $row13.Length; $i++) { $ws9.Cells.Item(13, $i + 1).Value = $row13[$i] }

# Row 14
$row14 = @(0.067,0.083,0.076,0.099,0.081,0.077,0.093,0.075,0.087,0.093,0.125,0.094)
for ($i = 0; $i -lt $row14.Length; $i++) { $ws9.Cells.Item(14, $i + 1).Value = $row14[$i] }

# Clear the stray highlight on G11 (it is excluded from the row-16 average)
$ws9.Range("G11").Interior.ColorIndex = -4142

# Row 16 - averages of rows 8:12 (G16 explicitly skips the G11 outlier)
$ws9.Range("B16").Formula = "=AVERAGE(B8:B12)"
$ws9.Range("C16").Formula = "=AVERAGE(C8:C12)"
$ws9.Range("D16").Formula = "=AVERAGE(D8:D12)"
$ws9.Range("E16").Formula = "=AVERAGE(E8:E12)"
$ws9.Range("F16").Formula = "=AVERAGE(F8:F12)"
$ws9.Range("G16").Formula = "=AVERAGE(G8:G10,G12)"
$ws9.Range("H16").Formula = "=AVERAGE(H8:H12)"
$ws9.Range("I16").Formula = "=AVERAGE(I8:I12)"
$ws9.Range("J16").Formula = "=AVERAGE(J8:J12)"
$ws9.Range("K16").Formula = "=AVERAGE(K8:K12)"

# Rows 18-27 - normalised-to-day-1 percentages, copy the blue fill style
# used on column C throughout the other sheets.
$ws8.Range("C19:C27").Copy()
$ws9.Range("C18:C27").PasteSpecial(-4122)

$ws9.Range("B18").Value = 1.3328
$ws9.Range("C18").Value = 100

$ws9.Range("B19").Value = 1.2154
$ws9.Range("C19").Formula = "=B19/1.3328*100"

$ws9.Range("B20").Value = 1.101
$ws9.Range("C20").Formula = "=B20/1.3328*100"

$ws9.Range("B21").Value = 0.93740000000000001
$ws9.Range("C21").Formula = "=B21/1.3328*100"

$ws9.Range("B22").Value = 0.85599999999999987
$ws9.Range("C22").Formula = "=B22/1.3328*100"

$ws9.Range("B23").Value = 0.64175000000000004
$ws9.Range("C23").Formula = "=B23/1.3328*100"

$ws9.Range("B24").Value = 0.55700000000000005
$ws9.Range("C24").Formula = "=B24/1.3328*100"

$ws9.Range("B25").Value = 0.43499999999999994
$ws9.Range("C25").Formula = "=B25/1.3328*100"

$ws9.Range("B26").Value = 0.30959999999999999
$ws9.Range("C26").Formula = "=B26/1.3328*100"

$ws9.Range("B27").Value = 0.32539999999999997
$ws9.Range("C27").Formula = "=B27/1.3328*100"

# Final selection/activation state for the new sheet matches the commit.
$ws9.Activate()
$ws9.Range("C18:C27").Select()
